$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (DAMSLTag, DialogAct) updates re-annotated by SGNN
$updates = @(
    @{ Row = 14; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 19; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 22; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 32; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 34; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 38; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 50; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 81; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 84; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 87; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 97; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 100; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 134; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 150; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 166; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 169; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 177; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 187; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 198; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 199; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 204; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 205; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 208; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 220; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 224; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 241; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 253; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 264; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 267; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 276; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 279; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 280; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 291; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 298; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 305; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 306; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 313; DAMSLTag = "ba"; DialogAct = "Appreciation" }
    @{ Row = 314; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
    @{ Row = 319; DAMSLTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 321; DAMSLTag = "b"; DialogAct = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()
